# This edit reorders/updates the data rows (rows 2-25) of the worksheet.
# Each destination row's editable columns (D and K:T) are replaced with the
# values that originally belonged to a different row (a permutation of the
# weekly records), while columns A:C and E:J (identical across all rows)
# remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values taken from the ORIGINAL workbook state)
$rowMap = @{
    2  = 8
    3  = 9
    4  = 13
    5  = 14
    6  = 10
    7  = 11
    8  = 22
    9  = 23
    10 = 4
    11 = 24
    12 = 25
    13 = 20
    14 = 21
    15 = 12
    16 = 18
    17 = 19
    18 = 15
    19 = 16
    20 = 17
    21 = 2
    22 = 3
    23 = 5
    24 = 6
    25 = 7
}

# Columns whose values vary per row and need to be taken from the mapped source row.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the original values for every row/column involved, before any writes happen,
# since this is a permutation (writes to one row must not affect reads for another).
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Apply the new values according to the mapping.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
